$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.252895951271057
$ws.Range("B1").Value = 1.834372878074646
$ws.Range("C1").Value = 4.172325611114502
$ws.Range("D1").Value = 3.209892988204956
$ws.Range("E1").Value = 1.161930680274963
